# Simplify database schema migration
# Append a new row (76) of data to each of the four log sheets, carrying
# forward the same layout/format used by the existing rows.

$wb = $excel.ActiveWorkbook

$rowData = @{
    1 = @{
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x3C"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H = 316
        I = 7
    }
    2 = @{
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x3C"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H = 316
        I = 25
    }
    3 = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x61"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H = 97
        I = 15
    }
    4 = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x77"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H = 119
        I = 9
    }
}

$newRow = 76
$timeValue = 45862.46436342593

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowData[$i]

    # Column A: timestamp, same date/time number format as the rows above it
    $ws.Cells.Item($newRow, 1).Value = $timeValue
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Columns B-E: hex byte strings stored as text
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    # Columns F-I: decimal counterparts stored as numbers
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = [double]$data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
